$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header for column C: "dummy_model_response" -> "model_1"
$ws.Range("C1").Value = "model_1"

# Update selection to C2 (matches the post-edit view state in the saved file)
$ws.Range("C2").Select()
